# repull data, push all data, mean calculation
# Update the dSF column (F) values with freshly-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = 6
$ws.Range("F9").Value = -7
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = -2
